# Working_time.xlsx -- "Add files via upload" re-save
#
# The sheet "額外列印" (extra printing) had its B3:B6 option labels
# retyped from the "塊" (piece) wording to the "種" (kind) wording, to
# match the vocabulary used on the "列印不同厚度" sheet. Both sheets keep
# their numeric Time column untouched.
$wb = $excel.ActiveWorkbook

$wsExtra = $wb.Worksheets.Item("額外列印")
$wsExtra.Range("B3").Value = "1種"
$wsExtra.Range("B4").Value = "2種"
$wsExtra.Range("B5").Value = "3種"
$wsExtra.Range("B6").Value = ">3種"

# Leave a cursor parked on the edited sheet ...
$wsExtra.Range("M16").Select()

# ... then hop back to "列印不同厚度", which was (and stays) the tab that
# was showing on screen when the file was saved, and park its cursor on
# the multi-cell range that was last highlighted there.
$wsThickness = $wb.Worksheets.Item("列印不同厚度")
$wsThickness.Range("L27:L29").Select()
